$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "날짜"
$ws.Range("B1").Value = "종목명"
$ws.Range("C1").Value = "티커"
$ws.Range("D1").Value = "종가"
$ws.Range("E1").Value = "RSI"
$ws.Range("F1").Value = "5일수익률"
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"
$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# give the newly added header cells (L1:O1) the same bold/centered/
# bordered look as the rest of the header row, without disturbing the
# text that was just written into them
$ws.Range("A1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows -------------------------------------------------------
# The "날짜" column holds a literal text value "2025-11-29", not a real
# date. Assigning that string straight to .Value makes Excel "helpfully"
# reinterpret it as a date serial number, so instead we stage it as a
# formula result in a scratch cell and paste-special just the value over
# - this keeps it a genuine text cell (t="s"), matching the source data.
$ws.Range("Z100").Formula = '="2025-11-29"'
$ws.Range("Z100").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("Z100").Clear()

# Row 2: SamsungHvyInd
$ws.Range("B2").Value = "SamsungHvyInd"
$ws.Range("C2").Value = "010140.KS"
$ws.Range("D2").Value = 24600
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = -2.96
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 60.6
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.36763896678245
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# Row 3: Hanwha Ocean
$ws.Range("Z100").Formula = '="2025-11-29"'
$ws.Range("Z100").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$ws.Range("Z100").Clear()
$ws.Range("B3").Value = "Hanwha Ocean"
$ws.Range("C3").Value = "042660.KS"
$ws.Range("D3").Value = 107800
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = -10.02
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 54.8
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.36763896678245
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# Row 4: HDKSOE
$ws.Range("Z100").Formula = '="2025-11-29"'
$ws.Range("Z100").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("Z100").Clear()
$ws.Range("B4").Value = "HDKSOE"
$ws.Range("C4").Value = "009540.KS"
$ws.Range("D4").Value = 410000
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = -2.38
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 70
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 51
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 85.36763896678245
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"

# Row 5: HD HYUNDAI MIPO
$ws.Range("Z100").Formula = '="2025-11-29"'
$ws.Range("Z100").Copy()
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("Z100").Clear()
$ws.Range("B5").Value = "HD HYUNDAI MIPO"
$ws.Range("C5").Value = "010620.KS"
$ws.Range("D5").Value = 223000
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = -4.09
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 46
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 48.6
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 85.36763896678245
$ws.Range("O5").Value = "🟢 완화적 (상승 우위)"
